$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 125, pushing existing rows 125-248 down to 126-249.
$ws.Rows.Item(125).Insert()

# Fill the newly inserted row 125 with data (copy constant template columns, set new values).
$ws.Cells.Item(125, 1).Value = 5
$ws.Cells.Item(125, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(125, 3).Value = "Maule"
$ws.Cells.Item(125, 4).Value = 44586
$ws.Cells.Item(125, 5).Value = 7
$ws.Cells.Item(125, 6).Value = 100114014
$ws.Cells.Item(125, 7).Value = "Betarraga"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 4000
$ws.Cells.Item(125, 11).Value = 700
$ws.Cells.Item(125, 12).Value = 700
$ws.Cells.Item(125, 13).Value = 700
$ws.Cells.Item(125, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(125, 15).Value = "Región del Maule"
$ws.Cells.Item(125, 16).Value = 140
$ws.Cells.Item(125, 17).Value = 5
$ws.Cells.Item(125, 18).Value = "Hortaliza"
